# Apply updates to difmedias_controles_staggered_variables_2015.xlsx
# - Update the dummy_d1 and cantidad_d1 rows with refreshed statistics
# - Add a new "ingreso" variable row (with N / Mean / SE / mean-difference values)
# - Move the footnote down one row and append "ingreso" to the recorded user command

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that receive plain numeric-looking text (e.g. "0.028", "860") must be
# forced to Text format first, otherwise Excel would store them as numbers.
$textCells = @("C6","E6","C12","E12","B16","C16","D16","E16","F16")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Updated statistics for dummy_d1 (row 6/7) ---
$ws.Range("C6").Value = "0.028"
$ws.Range("E6").Value = "0.111"
$ws.Range("G6").Value = "0.083***"
$ws.Range("C7").Value = "(0.006)"
$ws.Range("E7").Value = "(0.053)"

# --- Updated statistics for cantidad_d1 (row 12/13) ---
$ws.Range("C12").Value = "0.030"
$ws.Range("E12").Value = "0.111"
$ws.Range("G12").Value = "0.081**"
$ws.Range("C13").Value = "(0.006)"
$ws.Range("E13").Value = "(0.053)"

# --- New "ingreso" variable, replacing the old footnote row (row 16/17) ---
$ws.Range("A16").Value = "ingreso"
$ws.Range("B16").Value = "860"
$ws.Range("C16").Value = "2.088"
$ws.Range("D16").Value = "36"
$ws.Range("E16").Value = "2.830"
$ws.Range("F16").Value = "896"
$ws.Range("G16").Value = "0.742***"
$ws.Range("C17").Value = "(0.044)"
$ws.Range("E17").Value = "(0.086)"

# --- Footnote moves to row 18 and now mentions the new "ingreso" variable ---
$ws.Range("A18").Value = "If the table includes missing values (.n, .o, .v etc.) see the Missing values section in the help file for the Stata command iebaltab for definitions of these values. Significance: ***=.01, **=.05, *=.1. Full user input as written by user: [iebaltab dummy_jb dummy_d1 dummy_ara cantidad_jb cantidad_d1 cantidad_ara ingreso , groupvar(dummy_oxxo) control(0) savexlsx(difmedias_controles_staggered_variables_2015) replace] "

Write-Output "Edits applied"
